$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 new headers (Area / Atotal columns, plus a repeated Atotal/Qtotal pair) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- B2 / C2 become numeric 0 instead of "-" ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- New Area column formulas ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
# Fill G4:G15 together so the engine keeps them as one shared formula group
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

$ws.Range("H2").Formula = "=SUM(G2:G11)"

$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

$ws.Range("D2").Select() | Out-Null
